# Updated cryptos list on Fri Mar  1 16:46:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.478.01"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.393.44"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "404.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  -2.24%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.128"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.141"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "3.385.99"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "61.456.51"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000143"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.50%  "
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "82.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "310.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.07%  "
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "43.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "

# Row 38/39: Stacks and LidoDAOToken swap places with updated values
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "2.090.59"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("E51").Value = "  +16.66%  "
